$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 93.333336
$ws.Range("I6").Value = 93.333336
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 280.000008
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -168.000008
$ws.Range("N6").ClearContents()

$ws.Range("H18").Value = 1000000000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 1000000000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 1000000000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -1000000568

$ws.Range("H112").Value = 2267.875
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2267.875
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 6803.625
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -9019.625

$ws.Range("H132").Value = 5502.393
$ws.Range("I132").Value = 5691.3706
$ws.Range("J132").Value = 400
$ws.Range("K132").Value = 17074.1118
$ws.Range("L132").Value = 1200
$ws.Range("M132").Value = -14544.1118
$ws.Range("N132").Value = -6260

$ws.Range("H135").Value = 2350.889
$ws.Range("I135").Value = 1276.5
$ws.Range("J135").Value = 4499.6665
$ws.Range("K135").Value = 11488.5
$ws.Range("L135").Value = 40496.9985
$ws.Range("M135").Value = -8953.5

$ws.Range("H137").Value = 4548348
$ws.Range("I137").Value = 5558439.5
$ws.Range("J137").Value = 2936.5
$ws.Range("K137").Value = 16675318.5
$ws.Range("L137").Value = 8809.5
$ws.Range("M137").Value = -16672768.5
$ws.Range("N137").Value = -13909.5

$ws.Range("H138").Value = 5441
$ws.Range("I138").Value = 2331.5557
$ws.Range("J138").Value = 7306.6665
$ws.Range("K138").Value = 6994.6671
$ws.Range("L138").Value = 21919.9995
$ws.Range("M138").Value = -1854.6671
$ws.Range("N138").Value = -32199.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1389206.8
$ws.Range("I32").Value = 618709.3
$ws.Range("J32").Value = 47619050
$ws.Range("K32").Value = 618709.3
$ws.Range("L32").Value = 47619050
$ws.Range("M32").Value = -618422.3

$ws.Range("H61").Value = 2318.9092
$ws.Range("I61").Value = 2062.5652
$ws.Range("J61").Value = 2908.5
$ws.Range("K61").Value = 2062.5652
$ws.Range("L61").Value = 2908.5
$ws.Range("M61").Value = -1850.5652
$ws.Range("N61").Value = -3332.5

$ws.Range("H74").Value = 160434.17
$ws.Range("I74").Value = 214974.31
$ws.Range("J74").Value = 2873.7778
$ws.Range("K74").Value = 214974.31
$ws.Range("L74").Value = 2873.7778
$ws.Range("M74").Value = -214100.31

$ws.Range("H77").Value = 160434.17
$ws.Range("I77").Value = 214974.31
$ws.Range("J77").Value = 2873.7778
$ws.Range("K77").Value = 1074871.55
$ws.Range("L77").Value = 14368.889
$ws.Range("M77").Value = -1070503.55

$ws.Range("H97").Value = 1144.5
$ws.Range("I97").Value = 1113.5
$ws.Range("J97").Value = 1237.5
$ws.Range("K97").Value = 1113.5
$ws.Range("L97").Value = 1237.5
$ws.Range("M97").Value = -617.5

$ws.Range("H122").Value = 2944.0908
$ws.Range("I122").Value = 2592.7778
$ws.Range("J122").Value = 4525
$ws.Range("K122").Value = 7778.3334
$ws.Range("L122").Value = 13575
$ws.Range("M122").Value = -5328.3334

$ws.Range("H132").Value = 1587.5
$ws.Range("I132").Value = 1305
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3915
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1385
$ws.Range("N132").Value = -14060

$ws.Range("H136").Value = 2318.9092
$ws.Range("I136").Value = 2062.5652
$ws.Range("J136").Value = 2908.5
$ws.Range("K136").Value = 6187.6956
$ws.Range("L136").Value = 8725.5
$ws.Range("M136").Value = -3637.6956
$ws.Range("N136").Value = -13825.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 16252451
$ws.Range("I105").Value = 1252138.5
$ws.Range("J105").Value = 31252762
$ws.Range("K105").Value = 1252138.5
$ws.Range("L105").Value = 31252762
$ws.Range("M105").Value = -1250391.5

$ws.Range("H107").Value = 2332268.5
$ws.Range("I107").Value = 3206162.2
$ws.Range("J107").Value = 1885.1111
$ws.Range("K107").Value = 3206162.2
$ws.Range("L107").Value = 1885.1111
$ws.Range("M107").Value = -3204242.2
$ws.Range("N107").Value = -5725.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4319819
$ws.Range("I31").Value = 2479.25
$ws.Range("J31").Value = 7367353
$ws.Range("K31").Value = 2479.25
$ws.Range("L31").Value = 7367353
$ws.Range("M31").Value = -2184.25
$ws.Range("N31").Value = -7367943

$ws.Range("H34").Value = 4319819
$ws.Range("I34").Value = 2479.25
$ws.Range("J34").Value = 7367353
$ws.Range("K34").Value = 2479.25
$ws.Range("L34").Value = 7367353
$ws.Range("M34").Value = -2277.25
$ws.Range("N34").Value = -7367757

$ws.Range("H60").Value = 14649.5
$ws.Range("I60").Value = 9066.666999999999
$ws.Range("J60").Value = 17999.2
$ws.Range("K60").Value = 9066.666999999999
$ws.Range("L60").Value = 17999.2
$ws.Range("M60").Value = -8555.666999999999
$ws.Range("N60").Value = -19021.2

$ws.Range("H98").Value = 149999
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 149999
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 149999
$ws.Range("N98").Value = -154491

$ws.Range("H122").Value = 2452.6843
$ws.Range("I122").Value = 1879.4667
$ws.Range("J122").Value = 4602.25
$ws.Range("K122").Value = 5638.4001
$ws.Range("L122").Value = 13806.75
$ws.Range("M122").Value = -3188.4001
$ws.Range("N122").Value = -18706.75

$ws.Range("H132").Value = 3685.875
$ws.Range("I132").Value = 3476
$ws.Range("J132").Value = 5994.5
$ws.Range("K132").Value = 10428
$ws.Range("L132").Value = 17983.5
$ws.Range("M132").Value = -7898
$ws.Range("N132").Value = -23043.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1226.5
$ws.Range("I5").Value = 965.2
$ws.Range("J5").Value = 1662
$ws.Range("K5").Value = 2895.6
$ws.Range("L5").Value = 4986
$ws.Range("M5").Value = -2783.6
$ws.Range("N5").Value = -5210

$ws.Range("H9").Value = 37129756
$ws.Range("I9").Value = 55693988
$ws.Range("J9").Value = 1290.3334
$ws.Range("K9").Value = 167081964
$ws.Range("L9").Value = 3871.0002
$ws.Range("M9").Value = -167081740
$ws.Range("N9").Value = -4319.0002

$ws.Range("H34").Value = 469.66666
$ws.Range("I34").Value = 364
$ws.Range("J34").Value = 998
$ws.Range("K34").Value = 1092
$ws.Range("L34").Value = 2994
$ws.Range("M34").Value = -1008
$ws.Range("N34").Value = -3162

$ws.Range("H56").Value = 7580.75
$ws.Range("I56").Value = 7580.75
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 7580.75
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -7050.75

$ws.Range("H122").Value = 552.4
$ws.Range("I122").Value = 683.6667
$ws.Range("J122").Value = 496.14285
$ws.Range("K122").Value = 6153.0003
$ws.Range("L122").Value = 4465.28565
$ws.Range("M122").Value = -3703.0003
$ws.Range("N122").Value = -9365.28565

$ws.Range("H128").Value = 173327.67
$ws.Range("I128").Value = 173327.67
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 519983.01
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -515003.01

$ws.Range("H135").Value = 1226.5
$ws.Range("I135").Value = 965.2
$ws.Range("J135").Value = 1662
$ws.Range("K135").Value = 8686.800000000001
$ws.Range("L135").Value = 14958
$ws.Range("M135").Value = -6151.800000000001
$ws.Range("N135").Value = -20028

$ws.Range("H139").Value = 3509.7827
$ws.Range("I139").Value = 2807.4119
$ws.Range("J139").Value = 5499.8335
$ws.Range("K139").Value = 8422.235700000001
$ws.Range("L139").Value = 16499.5005
$ws.Range("M139").Value = -3282.235700000001
$ws.Range("N139").Value = -26779.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2247.4
$ws.Range("I2").Value = 2805.5
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 2805.5
$ws.Range("L2").Value = 15
$ws.Range("M2").Value = -2692.5
$ws.Range("N2").Value = -241

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H70").Value = 6408.6
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 6760.75
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 6760.75
$ws.Range("M70").Value = -4730

$ws.Range("H73").Value = 6408.6
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 6760.75
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 6760.75
$ws.Range("M73").Value = -4064

$ws.Range("H99").Value = 1588.5
$ws.Range("I99").Value = 1588.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1588.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 657.5

$ws.Range("H113").Value = 2283.5
$ws.Range("I113").Value = 2301
$ws.Range("J113").Value = 2196
$ws.Range("K113").Value = 2301
$ws.Range("L113").Value = 2196
$ws.Range("M113").Value = -131

$ws.Range("H122").Value = 2851969.2
$ws.Range("I122").Value = 5130476.5
$ws.Range("J122").Value = 3834.9167
$ws.Range("K122").Value = 15391429.5
$ws.Range("L122").Value = 11504.7501
$ws.Range("M122").Value = -15388979.5

$ws.Range("H126").Value = 6033
$ws.Range("I126").Value = 2617.5
$ws.Range("J126").Value = 9936.429
$ws.Range("K126").Value = 7852.5
$ws.Range("L126").Value = 29809.287
$ws.Range("M126").Value = -5382.5

$ws.Range("H132").Value = 2126.7273
$ws.Range("I132").Value = 1987.5714
$ws.Range("J132").Value = 2370.25
$ws.Range("K132").Value = 5962.7142
$ws.Range("L132").Value = 7110.75
$ws.Range("M132").Value = -3432.7142
$ws.Range("N132").Value = -12170.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3748.5
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 4664.6665
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 4664.6665
$ws.Range("M46").Value = -812

$ws.Range("H132").Value = 13073.917
$ws.Range("I132").Value = 7148
$ws.Range("J132").Value = 18999.834
$ws.Range("K132").Value = 21444
$ws.Range("L132").Value = 56999.50199999999
$ws.Range("M132").Value = -18914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H113").Value = 460.33334
$ws.Range("I113").Value = 439.6154
$ws.Range("J113").Value = 595
$ws.Range("K113").Value = 1318.8462
$ws.Range("L113").Value = 1785
$ws.Range("M113").Value = 851.1538
$ws.Range("N113").Value = -6125

$ws.Range("H132").Value = 2999.5557
$ws.Range("I132").Value = 2999.5557
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8998.667099999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6468.667099999999
